$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mediciones table updates (rows 5-17)
$ws.Range("C7").Value = 20
$ws.Range("F8").Value = 3
$ws.Range("D9").Value = 1
$ws.Range("G9").Value = 23

# More measurements
$ws.Range("C15").Value = 2000
$ws.Range("D16").Value = 3100
$ws.Range("G16").Value = 1

# Capataz name changes
$ws.Range("I17").Value = "PEPE"
$ws.Range("I16").Value = "CAMBIO"

# Observaciones notes
$ws.Range("J13").Value = "CAMBIO"
$ws.Range("J12").Value = "HOLA"

# Stray formatting on the cell right after the table (underline, no value)
$ws.Range("K17").Font.Underline = $true

# Selection state
$null = $ws.Range("R17").Select()
